$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update title / site text values
$ws.Range("A1").Value = "SmartFire"
$ws.Range("A2").Value = "smartfire.com.br"
$ws.Range("A5").Value = "Teste"

# 2. Update row 5 data values (B5 becomes 5, C5 becomes "Não encontrado")
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "Não encontrado"

# 3. Rebuild hyperlinks collection: keep only the one on A5 (same target
#    as before), drop the one on A6 (row 6 will be removed below).
$origUrl = "https://www.google.com/search?q=Cilindros%20Hidr%C3%A1ulicos%20De%20Alta%20Press%C3%A3o&start=0"
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A5"), $origUrl) | Out-Null
$ws.Range("A5").Style = "Normal"

# 4. Remove row 6 entirely (its data + hyperlink)
$ws.Rows.Item(6).Delete()

Write-Output "done"
